$d = $word.ActiveDocument

# Locate the "Pre-approvals for promotion" bullet paragraph by its text
# instead of a hard-coded index, so the script is robust to minor
# structural differences.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r")
    if ($ptext -eq "Pre-approvals for promotion") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# --- Step 1: split the "Pre-approvals for promotion" run into
#     "Pre-approvals for " + "promotion" (identical rPr on both) ---
$paraText = $p.Range.Text
$splitWord = "promotion"
$idx = $paraText.IndexOf($splitWord)
$start = $p.Range.Start + $idx
$end = $start + $splitWord.Length
$promoRange = $d.Range($start, $end)
# Toggling Bold on/off forces Word to re-serialize the paragraph text as
# two runs split at this boundary, without altering visible formatting.
$promoRange.Bold = 1
$promoRange.Bold = 0

# --- Step 2: add three new sub-bullets after the "Pre-approvals for
#     promotion" paragraph, inheriting its list level/indentation/style ---
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($targetIndex + 1)
$p1.Range.Text = "Market segmentation: assign risk bands according to prob of degfault. Then apply loan terms accordingly. "

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($targetIndex + 2)
$p2.Range.Text = "Optimize reserves"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($targetIndex + 3)
$p3.Range.Text = "Re-Calibrate according to risk apetite. "
